# Updated cryptos list with GitHub Actions
# Applies the latest price / 1h-volume-change snapshot to Sheet1.
# Each target cell stores its value as text (inline/shared string), so we
# force the cell to Text format before assigning, then restore the default
# "Normal" style afterward so no stray number-format style lingers on the
# cell (matches the original workbook, which has no explicit style on
# these cells).
function Set-CellText($ws, $ref, $text) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws 'D2' '65.704.41'
Set-CellText $ws 'E2' '  +0.88%  '
Set-CellText $ws 'D3' '3.579.86'
Set-CellText $ws 'E3' '  +1.49%  '
Set-CellText $ws 'E4' '  +0.00%  '
Set-CellText $ws 'D5' '601.88'
Set-CellText $ws 'E5' '  +0.84%  '
Set-CellText $ws 'D6' '137.67'
Set-CellText $ws 'E6' '  -0.30%  '
Set-CellText $ws 'D7' '3.579.33'
Set-CellText $ws 'E7' '  +1.45%  '
Set-CellText $ws 'E8' '  -0.02%  '
Set-CellText $ws 'D9' '0.496'
Set-CellText $ws 'E9' '  +0.64%  '
Set-CellText $ws 'E10' '  +0.72%  '
Set-CellText $ws 'E11' '  +5.44%  '
Set-CellText $ws 'E12' '  +1.02%  '
Set-CellText $ws 'D13' '4.188.63'
Set-CellText $ws 'E13' '  +1.53%  '
Set-CellText $ws 'E14' '  +3.78%  '
Set-CellText $ws 'E15' '  +0.98%  '
Set-CellText $ws 'D16' '3.583.11'
Set-CellText $ws 'E16' '  +1.65%  '
Set-CellText $ws 'E17' '  -0.16%  '
Set-CellText $ws 'D18' '65.757.38'
Set-CellText $ws 'E18' '  +0.83%  '
Set-CellText $ws 'D19' '10.00'
Set-CellText $ws 'E19' '  -2.85%  '
Set-CellText $ws 'D20' '14.62'
Set-CellText $ws 'E20' '  +2.47%  '
Set-CellText $ws 'E21' '  -1.18%  '
Set-CellText $ws 'D22' '395.02'
Set-CellText $ws 'E22' '  +0.75%  '
Set-CellText $ws 'E23' '  +3.41%  '
Set-CellText $ws 'D24' '3.725.26'
Set-CellText $ws 'E24' '  +1.56%  '
Set-CellText $ws 'D25' '74.15'
Set-CellText $ws 'E25' '  +0.47%  '
Set-CellText $ws 'E26' '  +0.00%  '
Set-CellText $ws 'E27' '  +2.95%  '
Set-CellText $ws 'E28' '  +5.48%  '
Set-CellText $ws 'D29' '1.63'
Set-CellText $ws 'E29' '  +29.89%  '
Set-CellText $ws 'D30' '2.38'
Set-CellText $ws 'E30' '  +3.98%  '
Set-CellText $ws 'E31' '  +4.98%  '
Set-CellText $ws 'E32' '  -0.29%  '
Set-CellText $ws 'D33' '3.584.17'
Set-CellText $ws 'E33' '  +1.26%  '
Set-CellText $ws 'D34' '24.50'
Set-CellText $ws 'E34' '  +2.97%  '
Set-CellText $ws 'E36' '  +2.28%  '
Set-CellText $ws 'D37' '5.39'
Set-CellText $ws 'E37' '  +8.46%  '
Set-CellText $ws 'E38' '  +5.18%  '
Set-CellText $ws 'E39' '  +1.87%  '
Set-CellText $ws 'D40' '168.88'
Set-CellText $ws 'E40' '  +0.16%  '
Set-CellText $ws 'D41' '0.0835'
Set-CellText $ws 'E41' '  +4.67%  '
Set-CellText $ws 'E42' '  +2.22%  '
Set-CellText $ws 'D43' '26.95'
Set-CellText $ws 'E43' '  +3.81%  '
Set-CellText $ws 'D44' '1.28'
Set-CellText $ws 'E44' '  +8.59%  '
Set-CellText $ws 'D45' '43.15'
Set-CellText $ws 'E45' '  +0.90%  '
Set-CellText $ws 'E46' '  +2.97%  '
Set-CellText $ws 'E47' '  -0.03%  '
Set-CellText $ws 'E48' '  +2.13%  '
Set-CellText $ws 'D49' '7.02'
Set-CellText $ws 'E49' '  +3.38%  '
Set-CellText $ws 'D50' '2.454.44'
Set-CellText $ws 'D51' '0.908'
Set-CellText $ws 'E51' '  +10.54%  '
